# Update the workbook so each struct attribute becomes a column header,
# and any attribute whose value list has more than one option gets a
# dropdown (data validation) on its column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace header row: Número/Mês/Ano -> Name/Age/Email, and add Games.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Games"

# Remove the old Month/Year dropdown validations from columns B and C.
$ws.Range("B2:B1000").Validation.Delete()
$ws.Range("C2:C1000").Validation.Delete()

# Add a new dropdown validation for the Games column.
$ws.Range("D2:D100000").Validation.Add(3, 1, 1, '"Super Mario,SONIC,Zelda,GTA"')
